# Updates cryptos list data (price + 1h volume change) to reflect the
# latest scrape, including the EthereumClassic/Hedera rank swap (rows 33/34).
# Numeric-looking price strings in column D are forced to stay text (they are
# stored as plain strings in the sheet, not real numbers) by briefly switching
# the cell to a text number format, then restoring the original Normal style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.763.16"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "3.073.22"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").Value = "3.072.95"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("E15").Value = "  +1.23%  "

$ws.Range("D16").Value = "3.581.38"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "63.682.44"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Value = "3.073.98"
$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").Value = "  -2.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.23%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.112"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.72%  "

$ws.Range("E35").Value = "  +1.61%  "

$ws.Range("D36").Value = "0.0₃0822"
$ws.Range("E36").Value = "  -3.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("E39").Value = "  -5.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.02%  "

$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("D46").Value = "2.842.54"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
